$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 664
    $ws.Range("F3").Value = 504
    $ws.Range("F5").Value = 20
    $ws.Range("F8").Value = 2605
    $ws.Range("F9").Value = 4142
}
